# Update the f0 (column K) values for rows 27-51 on Sheet1 from 2000 to 4000.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 27; $row -le 51; $row++) {
    $cell = $ws.Cells.Item($row, 11)
    if ($cell.Value2 -eq 2000) {
        $cell.Value = 4000
    }
}
